$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 202
$ws.Range("H2").Value = "bedrooms"
$ws.Range("I2").Value = "target"
$ws.Range("K2").Value = "j"
$ws.Range("L2").Value = "stimuli/img_uxxo0.png"
$ws.Range("M2").Value = 71.74418604651163
$ws.Range("N2").Value = 48.44186046511628
$ws.Range("O2").Value = 60.09302325581395
$ws.Range("P2").Value = 43
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 5

# Row 3
$ws.Range("F3").Value = 203
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "distractor"
$ws.Range("K3").Value = "f"
$ws.Range("L3").Value = "stimuli/img_di49w.png"
$ws.Range("M3").Value = 73.03125
$ws.Range("N3").Value = 52.25
$ws.Range("O3").Value = 62.640625
$ws.Range("P3").Value = 32
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 6
$ws.Range("S3").Value = 6

# Row 4
$ws.Range("F4").Value = 204
$ws.Range("H4").Value = "bedrooms"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_fqgem.png"
$ws.Range("M4").Value = 80.75
$ws.Range("N4").Value = 61.475
$ws.Range("O4").Value = 71.1125
$ws.Range("P4").Value = 40
$ws.Range("Q4").Value = 8
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = 8

# Row 5
$ws.Range("F5").Value = 205
$ws.Range("H5").Value = "bedrooms"
$ws.Range("I5").Value = "target"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_okvvw.png"
$ws.Range("M5").Value = 50.58333333333334
$ws.Range("N5").Value = 32.11111111111111
$ws.Range("O5").Value = 41.34722222222223
$ws.Range("P5").Value = 36
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2

# Row 6
$ws.Range("F6").Value = 206
$ws.Range("H6").Value = "living_rooms"
$ws.Range("I6").Value = "distractor"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_x4bln.png"
$ws.Range("M6").Value = 76.34042553191489
$ws.Range("N6").Value = 59.51063829787234
$ws.Range("O6").Value = 67.92553191489361
$ws.Range("P6").Value = 47
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 7

# Row 7
$ws.Range("F7").Value = 207
$ws.Range("H7").Value = "bedrooms"
$ws.Range("I7").Value = "target"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_71mhq.png"
$ws.Range("M7").Value = 69.34210526315789
$ws.Range("N7").Value = 47.02631578947368
$ws.Range("O7").Value = 58.18421052631579
$ws.Range("P7").Value = 38
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = 5
$ws.Range("S7").Value = 5

# Row 8
$ws.Range("F8").Value = 208
$ws.Range("H8").Value = "kitchens"
$ws.Range("I8").Value = "distractor"
$ws.Range("K8").Value = "f"
$ws.Range("L8").Value = "stimuli/img_cnyac.png"
$ws.Range("M8").Value = 69.14705882352941
$ws.Range("N8").Value = 47.8235294117647
$ws.Range("O8").Value = 58.48529411764706
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 5
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 5

# Row 9
$ws.Range("F9").Value = 209
$ws.Range("H9").Value = "kitchens"
$ws.Range("I9").Value = "distractor"
$ws.Range("K9").Value = "f"
$ws.Range("L9").Value = "stimuli/img_p3hpc.png"
$ws.Range("M9").Value = 72.83333333333333
$ws.Range("N9").Value = 52.22222222222222
$ws.Range("O9").Value = 62.52777777777777
$ws.Range("P9").Value = 36
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 6

# Row 10
$ws.Range("F10").Value = 210
$ws.Range("H10").Value = "bedrooms"
$ws.Range("I10").Value = "target"
$ws.Range("K10").Value = "j"
$ws.Range("L10").Value = "stimuli/img_x0u5z.png"
$ws.Range("M10").Value = 92
$ws.Range("N10").Value = 78.16216216216216
$ws.Range("O10").Value = 85.08108108108108
$ws.Range("P10").Value = 37
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = 10
$ws.Range("S10").Value = 10

# Row 11
$ws.Range("F11").Value = 211
$ws.Range("H11").Value = "bedrooms"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_a9acb.png"
$ws.Range("M11").Value = 77.11428571428571
$ws.Range("N11").Value = 58.42857142857143
$ws.Range("O11").Value = 67.77142857142857
$ws.Range("P11").Value = 35
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7

# Row 12
$ws.Range("F12").Value = 212
$ws.Range("H12").Value = "bedrooms"
$ws.Range("I12").Value = "target"
$ws.Range("K12").Value = "j"
$ws.Range("L12").Value = "stimuli/img_ybbmx.png"
$ws.Range("M12").Value = 55.24324324324324
$ws.Range("N12").Value = 36.75675675675676
$ws.Range("O12").Value = 46
$ws.Range("P12").Value = 37
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 3

# Row 13
$ws.Range("F13").Value = 213
$ws.Range("H13").Value = "bedrooms"
$ws.Range("I13").Value = "target"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_5p2ql.png"
$ws.Range("M13").Value = 89.19565217391305
$ws.Range("N13").Value = 72.52173913043478
$ws.Range("O13").Value = 80.85869565217391
$ws.Range("P13").Value = 46
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 10

# Row 14
$ws.Range("F14").Value = 214
$ws.Range("H14").Value = "bedrooms"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_2js6m.png"
$ws.Range("M14").Value = 40.02777777777778
$ws.Range("N14").Value = 20.88888888888889
$ws.Range("O14").Value = 30.45833333333334
$ws.Range("P14").Value = 36
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 2
$ws.Range("S14").Value = 2

# Row 15
$ws.Range("F15").Value = 215
$ws.Range("H15").Value = "bedrooms"
$ws.Range("I15").Value = "target"
$ws.Range("K15").Value = "j"
$ws.Range("L15").Value = "stimuli/img_rvssl.png"
$ws.Range("M15").Value = 74.25
$ws.Range("N15").Value = 54.33333333333334
$ws.Range("O15").Value = 64.29166666666667
$ws.Range("P15").Value = 36
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6

# Row 16
$ws.Range("F16").Value = 216
$ws.Range("H16").Value = "bedrooms"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_2pk6v.png"
$ws.Range("M16").Value = 85.08108108108108
$ws.Range("N16").Value = 66.16216216216216
$ws.Range("O16").Value = 75.62162162162161
$ws.Range("P16").Value = 37
$ws.Range("Q16").Value = 9
$ws.Range("R16").Value = 9
$ws.Range("S16").Value = 9

# Row 17
$ws.Range("F17").Value = 217
$ws.Range("H17").Value = "bedrooms"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_oou46.png"
$ws.Range("M17").Value = 75.70270270270271
$ws.Range("N17").Value = 54.86486486486486
$ws.Range("O17").Value = 65.28378378378379
$ws.Range("P17").Value = 37
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = 6
$ws.Range("S17").Value = 6

# Row 18
$ws.Range("F18").Value = 218
$ws.Range("H18").Value = "living_rooms"
$ws.Range("I18").Value = "distractor"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/img_z4jxm.png"
$ws.Range("M18").Value = 88.30952380952381
$ws.Range("N18").Value = 72.64285714285714
$ws.Range("O18").Value = 80.47619047619048
$ws.Range("P18").Value = 42
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = 10
$ws.Range("S18").Value = 10

# Row 19
$ws.Range("F19").Value = 219
$ws.Range("H19").Value = "bedrooms"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_bj2gr.png"
$ws.Range("M19").Value = 65.25
$ws.Range("N19").Value = 44.8
$ws.Range("O19").Value = 55.025
$ws.Range("P19").Value = 40
$ws.Range("Q19").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("S19").Value = 4

# Row 20
$ws.Range("F20").Value = 220
$ws.Range("H20").Value = "bedrooms"
$ws.Range("I20").Value = "target"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_v8dra.png"
$ws.Range("M20").Value = 61.77272727272727
$ws.Range("N20").Value = 38.79545454545455
$ws.Range("O20").Value = 50.28409090909091
$ws.Range("P20").Value = 44
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 3
$ws.Range("S20").Value = 3

# Row 21
$ws.Range("F21").Value = 221
$ws.Range("H21").Value = "bedrooms"
$ws.Range("I21").Value = "target"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_le8uf.png"
$ws.Range("M21").Value = 12.88888888888889
$ws.Range("N21").Value = 9.222222222222221
$ws.Range("O21").Value = 11.05555555555556
$ws.Range("P21").Value = 36
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1

# Row 22
$ws.Range("F22").Value = 222
$ws.Range("H22").Value = "bedrooms"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_h0hbk.png"
$ws.Range("M22").Value = 86.80952380952381
$ws.Range("N22").Value = 69.19047619047619
$ws.Range("O22").Value = 78
$ws.Range("P22").Value = 42
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9

# Row 23
$ws.Range("F23").Value = 223
$ws.Range("H23").Value = "bedrooms"
$ws.Range("I23").Value = "target"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_th7xh.png"
$ws.Range("M23").Value = 82.35897435897436
$ws.Range("N23").Value = 65.53846153846153
$ws.Range("O23").Value = 73.94871794871796
$ws.Range("P23").Value = 39
$ws.Range("Q23").Value = 8
$ws.Range("R23").Value = 8
$ws.Range("S23").Value = 8

# Row 24
$ws.Range("F24").Value = 224
$ws.Range("H24").Value = "bedrooms"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_qgbyn.png"
$ws.Range("M24").Value = 65.08108108108108
$ws.Range("N24").Value = 40.10810810810811
$ws.Range("O24").Value = 52.5945945945946
$ws.Range("P24").Value = 37
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 4
$ws.Range("S24").Value = 4

# Row 25
$ws.Range("F25").Value = 225
$ws.Range("H25").Value = "bedrooms"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_5m6x4.png"
$ws.Range("M25").Value = 80.23076923076923
$ws.Range("N25").Value = 58.41025641025641
$ws.Range("O25").Value = 69.32051282051282
$ws.Range("P25").Value = 39
$ws.Range("Q25").Value = 7
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 7

# Row 26
$ws.Range("F26").Value = 226
$ws.Range("H26").Value = "kitchens"
$ws.Range("I26").Value = "distractor"
$ws.Range("K26").Value = "f"
$ws.Range("L26").Value = "stimuli/img_2b8fp.png"
$ws.Range("M26").Value = 73.89189189189189
$ws.Range("N26").Value = 51.45945945945946
$ws.Range("O26").Value = 62.67567567567568
$ws.Range("P26").Value = 37
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = 6
$ws.Range("S26").Value = 6

# Row 27
$ws.Range("F27").Value = 227
$ws.Range("H27").Value = "bedrooms"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_wyctg.png"
$ws.Range("M27").Value = 33.44736842105263
$ws.Range("N27").Value = 11.39473684210526
$ws.Range("O27").Value = 22.42105263157895
$ws.Range("P27").Value = 38
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
